# Updated Beq and Jump
# This script edits the ISA table:
#  1. "Init R1, 1" row -> prepend a leading space to the rightmost "Example Coding" cell.
#  2. "Beq" row -> binary coding cells change from "10111" pattern to "1011" + an extra
#     cyan "i" (i.e. the immediate field grows from "ii" to "iii" while the fixed opcode
#     bits shrink from "10111" to "1011").
#  3. "J" (jump) row -> same kind of change: opcode "10001"/"10001" -> "1000"/"1000" and
#     immediate field "ii"/"11" -> "iii"/"111".

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wdTurquoise = 3

function Insert-LeadingSpace($row, $col) {
    $cell = $t.Cell($row, $col)
    $full = $cell.Range
    $full.InsertBefore(" ")

    $t2 = $d.Tables.Item(1)
    $cell2 = $t2.Cell($row, $col)
    $full2 = $cell2.Range
    $newChar = $d.Range($full2.Start, $full2.Start + 1)
    $newChar.Font.Size = 15
}

function Insert-CyanCharBefore($row, $col, $offset, [string]$ch) {
    # Inserts $ch (plain text run) right before character offset $offset (0-based,
    # relative to the start of the cell) inside the cell at (row, col), then colors
    # the newly inserted character cyan (turquoise highlight) and restores its font size.
    $cell = $t.Cell($row, $col)
    $full = $cell.Range
    $pos = $full.Start + $offset
    $insertRange = $d.Range($pos, $pos)
    $insertRange.InsertBefore($ch)

    $t2 = $d.Tables.Item(1)
    $cell2 = $t2.Cell($row, $col)
    $full2 = $cell2.Range
    $newChar = $d.Range($full2.Start + $offset, $full2.Start + $offset + $ch.Length)
    $newChar.Font.Size = 15
    $newChar.Font.HighlightColorIndex = $wdTurquoise
}

# ---------------------------------------------------------------------------
# 1. "Init R1, 1" row (row 2): last column "0 000 01 01" -> " 0 000 01 01"
# ---------------------------------------------------------------------------
Insert-LeadingSpace 2 6

# ---------------------------------------------------------------------------
# 2. "Beq" row (row 10)
# ---------------------------------------------------------------------------

# Column 4 (binary for "Beq imm"): "10111 ii" -> "1011 iii"
$cell = $t.Cell(10, 4)
$cell.Range.Find.Execute("10111 ", $true, $false, $false, $false, $false, $true, 1, $false, "1011 ", 1) | Out-Null
# Now text is "1011 ii" ; insert a cyan "i" right before the existing cyan "ii" (offset 5)
Insert-CyanCharBefore 10 4 5 "i"

# Column 6 (binary for "Beq 4"): "0 10111 11" -> " 0 1011 111"
Insert-LeadingSpace 10 6
$cell = $t.Cell(10, 6)
$cell.Range.Find.Execute("10111", $true, $false, $false, $false, $false, $true, 1, $false, "1011", 1) | Out-Null
$cell = $t.Cell(10, 6)
$cell.Range.Find.Execute("11", $true, $false, $false, $false, $false, $true, 1, $false, "111", 1) | Out-Null

# ---------------------------------------------------------------------------
# 3. "J" row (row 11)
# ---------------------------------------------------------------------------

# Column 4 (binary for "J imm"): "10001 ii" -> "1000 iii"
$cell = $t.Cell(11, 4)
$cell.Range.Find.Execute("10001 ", $true, $false, $false, $false, $false, $true, 1, $false, "1000 ", 1) | Out-Null
# Now text is "1000 ii" ; insert a cyan "i" right before the existing cyan "ii" (offset 5)
Insert-CyanCharBefore 11 4 5 "i"

# Column 6 (binary for "J 3"): " 0 10001 11" -> " 0 1000 111"
$cell = $t.Cell(11, 6)
$cell.Range.Find.Execute("10001", $true, $false, $false, $false, $false, $true, 1, $false, "1000", 1) | Out-Null
$cell = $t.Cell(11, 6)
$cell.Range.Find.Execute("11", $true, $false, $false, $false, $false, $true, 1, $false, "111", 1) | Out-Null
